$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "Assignment Due" -> "Assignment", add new "Exam" column (F) ---
$ws.Range("E1").Value = "Assignment"

# F1 needs the same "Times New Roman" header style as the other header cells (E1/D1/etc).
# Copy formatting from D1 (already styled) onto F1, then set its value.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Exam"

# --- Week 3 (row 5): remove "Exam 1" text, add Assignment #1 marker in new column F ---
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 1

# --- Week 4 (row 6): remove "Privacy Statement (Initial)" ---
$ws.Range("E6").Value = ""

# --- Week 7 (row 9): unit label removed, topic/assignment updated ---
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "Color, stacking and dodging bars, and faceting"
$ws.Range("E9").Value = "Data ethics statement (draft)"

# --- Week 8 (row 10): topic/assignment updated ---
$ws.Range("C10").Value = "Summarizing data in tables"
$ws.Range("E10").Value = "Data ethics statement feedback"

# --- Week 9 (row 11): topic updated, new Assignment #2 marker added to F ---
$ws.Range("C11").Value = "More on tables and a gentle introduction to modeling"
$ws.Range("D11").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("F11").Value = 2

# --- Week 10 (row 12): unit label updated, exam removed entirely (incl. formatting) ---
$ws.Range("B12").Value = "3: Drawing conclusions"
$ws.Range("E12").Clear()

# --- Week 12 (row 14): blank styled placeholder added in column E ---
$ws.Range("D14").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# --- Week 13 (row 15): new Assignment #3 marker added to F ---
$ws.Range("F15").Value = 3

# --- Week 14 (row 16): final note text updated ---
$ws.Range("E16").Value = "Revised data ethics statement; Feedback on final presentations"

# --- Misc view state ---
$ws.Range("E17").Select()
